# Commit: "popup cambio modalita e clutch ok"
#
# Rewrites the Translation sheet's text table (Table8, columns B:F =
# Text ID / Typography Name / Alignment / GB / Direction) for rows 236-262.
# The 8 'Extra' placeholder rows that used to occupy rows 236-243 are
# collapsed to 2 rows (236-237); the existing 'LittleMedium' rows shift up
# to 238-251; 5 new 'Medium' race-mode labels are inserted at 252-256
# (ACCELERATION / ENDURANCE / AUTOCROSS / SKIDPAD / <value>); and a further
# 6 new 'Extra' popup rows are appended at 257-262.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    ,@(236, "SingleUseId253", "Extra", "Center", "<value>", "LTR")
    ,@(237, "SingleUseId254", "Extra", "Left", "N", "LTR")
    ,@(238, "SingleUseId258", "LittleMedium", "Center", "<value>", "LTR")
    ,@(239, "SingleUseId259", "LittleMedium", "Left", "ind_5", "LTR")
    ,@(240, "SingleUseId260", "LittleMedium", "Center", "<value>", "LTR")
    ,@(241, "SingleUseId261", "LittleMedium", "Left", "0.00", "LTR")
    ,@(242, "SingleUseId262", "LittleMedium", "Center", "<value>", "LTR")
    ,@(243, "SingleUseId263", "LittleMedium", "Left", "0.00", "LTR")
    ,@(244, "SingleUseId264", "LittleMedium", "Center", "<value>", "LTR")
    ,@(245, "SingleUseId265", "LittleMedium", "Left", "ind_5", "LTR")
    ,@(246, "SingleUseId266", "LittleMedium", "Center", "<value>", "LTR")
    ,@(247, "SingleUseId267", "LittleMedium", "Left", "0.00", "LTR")
    ,@(248, "SingleUseId268", "LittleMedium", "Center", "<value>", "LTR")
    ,@(249, "SingleUseId269", "LittleMedium", "Left", "ind_5", "LTR")
    ,@(250, "SingleUseId270", "LittleMedium", "Center", "<value>", "LTR")
    ,@(251, "SingleUseId271", "LittleMedium", "Left", "0.00", "LTR")
    ,@(252, "SingleUseId272", "Medium", "Center", "ACCELERATION", "LTR")
    ,@(253, "SingleUseId273", "Medium", "Center", "ENDURANCE", "LTR")
    ,@(254, "SingleUseId274", "Medium", "Center", "AUTOCROSS", "LTR")
    ,@(255, "SingleUseId275", "Medium", "Center", "SKIDPAD", "LTR")
    ,@(256, "SingleUseId276", "Medium", "Center", "<value>", "LTR")
    ,@(257, "SingleUseId277", "Extra", "Center", "<value>", "LTR")
    ,@(258, "SingleUseId278", "Extra", "Left", "N", "LTR")
    ,@(259, "SingleUseId279", "Extra", "Center", "<value>", "LTR")
    ,@(260, "SingleUseId280", "Extra", "Left", "N", "LTR")
    ,@(261, "SingleUseId281", "Extra", "Center", "<value>", "LTR")
    ,@(262, "SingleUseId282", "Extra", "Left", "N", "LTR")
)

# Helper: write a literal text value into a cell without letting Excel
# auto-coerce numeric-looking strings (e.g. "0.00") into numbers, and
# without leaving any residual custom cell style behind.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

foreach ($row in $rows) {
    $r = $row[0]
    Set-TextValue $ws.Cells.Item($r, 2) $row[1]   # B: Text ID
    Set-TextValue $ws.Cells.Item($r, 3) $row[2]   # C: Typography Name
    Set-TextValue $ws.Cells.Item($r, 4) $row[3]   # D: Alignment
    Set-TextValue $ws.Cells.Item($r, 5) $row[4]   # E: GB (text value)
    Set-TextValue $ws.Cells.Item($r, 6) $row[5]   # F: Direction
}
